$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.411890494488538
$ws.Range("C2").Value = 9.921937215030141
$ws.Range("D2").Value = -8.223290560045655
$ws.Range("E2").Value = -0.07183182115535164
$ws.Range("F2").Value = 2.081653799935187
$ws.Range("G2").Value = -1.355506186156552
$ws.Range("H2").Value = -1.41527269825023
$ws.Range("I2").Value = 0.6268712068528354
$ws.Range("J2").Value = 0.04461279023355938
$ws.Range("K2").Value = 0.1404525618062259
$ws.Range("B3").Value = 9.006385917080065
$ws.Range("C3").Value = -8.778729474695623
$ws.Range("D3").Value = -0.3907515313531172
$ws.Range("E3").Value = 1.872451415467592
$ws.Range("F3").Value = -1.50872361796536
$ws.Range("G3").Value = -1.540727567856255
$ws.Range("H3").Value = 0.5152988150593795
$ws.Range("I3").Value = -0.06003496370037081
$ws.Range("J3").Value = 0.03926128324195249
$ws.Range("K3").Value = -0.219352408342644
$ws.Range("B4").Value = -15.94395935728257
$ws.Range("C4").Value = -6.165406328900056
$ws.Range("D4").Value = -2.769606545464427
$ws.Range("E4").Value = -5.259251446597505
$ws.Range("F4").Value = -4.596812736670638
$ws.Range("G4").Value = -2.004405956105766
$ws.Range("H4").Value = -2.16751307464721
$ws.Range("I4").Value = -1.752444458675451
$ws.Range("J4").Value = -1.769680916244841
$ws.Range("K4").Value = -0.749069286040335
$ws.Range("B5").Value = -4.485225182330832
$ws.Range("C5").Value = 4.646438217044278
$ws.Range("D5").Value = -3.011195130372411
$ws.Range("E5").Value = -0.6154517214814641
$ws.Range("F5").Value = -0.0262240541957911
$ws.Range("G5").Value = 0.2444931714231107
$ws.Range("H5").Value = -0.1624340801744875
$ws.Range("I5").Value = -0.1271083699350421
$ws.Range("J5").Value = 0.5342533536669872
$ws.Range("K5").Value = 0.5366457149346298
$ws.Range("B6").Value = 0.9192316903262832
$ws.Range("C6").Value = -1.148324475608788
$ws.Range("D6").Value = -1.376360645652444
$ws.Range("E6").Value = 0.5240069414264609
$ws.Range("F6").Value = -0.002760069500128914
$ws.Range("G6").Value = 0.09815006044653379
$ws.Range("H6").Value = -0.1708893901647303
$ws.Range("I6").Value = 0.6659408511361884
$ws.Range("J6").Value = 0.5666756367783218
$ws.Range("K6").Value = 0.2125075656625323
$ws.Range("B7").Value = -0.740092021987018
$ws.Range("C7").Value = -1.12473760176605
$ws.Range("D7").Value = 0.4421008876033256
$ws.Range("E7").Value = 0.07233197988115672
$ws.Range("F7").Value = 0.2040572209602576
$ws.Range("G7").Value = -0.1124727440060749
$ws.Range("H7").Value = 0.7323795349662799
$ws.Range("I7").Value = 0.6413575953891584
$ws.Range("J7").Value = 0.2828208575635111
$ws.Range("K7").Value = 0.4305195042266514
$ws.Range("B8").Value = -1.230236813614865
$ws.Range("C8").Value = 0.3159956036379677
$ws.Range("D8").Value = 0.2405800485143884
$ws.Range("E8").Value = 0.2647213135641399
$ws.Range("F8").Value = -0.101339904631442
$ws.Range("G8").Value = 0.7952422862932105
$ws.Range("H8").Value = 0.6987091231829015
$ws.Range("I8").Value = 0.3273260345678901
$ws.Range("J8").Value = 0.48159337525307
$ws.Range("K8").Value = 0.6412619431822899
$ws.Range("B9").Value = -0.2884450162097709
$ws.Range("C9").Value = 0.176376671159695
$ws.Range("D9").Value = 0.5562686216213567
$ws.Range("E9").Value = -0.1390467807173756
$ws.Range("F9").Value = 0.7866582662197494
$ws.Range("G9").Value = 0.7824502007802698
$ws.Range("H9").Value = 0.3621531794959351
$ws.Range("I9").Value = 0.508157758305367
$ws.Range("J9").Value = 0.6867607839288887
$ws.Range("K9").Value = 0.1781548843816368
$ws.Range("B10").Value = 0.08873584232078371
$ws.Range("C10").Value = 0.5038494195916882
$ws.Range("D10").Value = -0.08040038373236996
$ws.Range("E10").Value = 0.7915874115568508
$ws.Range("F10").Value = 0.7713569177313746
$ws.Range("G10").Value = 0.3749198787210216
$ws.Range("H10").Value = 0.5168873173527261
$ws.Range("I10").Value = 0.6892627280777406
$ws.Range("J10").Value = 0.1844308218533179
$ws.Range("K10").Value = 0.4621062544735015
$ws.Range("B11").Value = 0.5262773992246967
$ws.Range("C11").Value = -0.07779684416992128
$ws.Range("D11").Value = 0.7695197259028355
$ws.Range("E11").Value = 0.7653749374270196
$ws.Range("F11").Value = 0.3695766285386105
$ws.Range("G11").Value = 0.5065187839823598
$ws.Range("H11").Value = 0.6807854059541167
$ws.Range("I11").Value = 0.1766585135465071
$ws.Range("J11").Value = 0.4534822809405329
$ws.Range("K11").Value = 0.2839429537672726
$ws.Range("B12").Value = -0.05637227271270118
$ws.Range("C12").Value = 0.8896975747020479
$ws.Range("D12").Value = 0.6761805930132365
$ws.Range("E12").Value = 0.3398496246900327
$ws.Range("F12").Value = 0.51514490401227
$ws.Range("G12").Value = 0.6548452325286815
$ws.Range("H12").Value = 0.1538799505182796
$ws.Range("I12").Value = 0.4398642457919151
$ws.Range("J12").Value = 0.2654221565820968
$ws.Range("K12").Value = 0.5706702220727796
$ws.Range("B13").Value = 0.8473769423816584
$ws.Range("C13").Value = 0.6469761191575033
$ws.Range("D13").Value = 0.3368290248851115
$ws.Range("E13").Value = 0.4962836204007895
$ws.Range("F13").Value = 0.6345141014634773
$ws.Range("G13").Value = 0.1393497950840318
$ws.Range("H13").Value = 0.4232232002996471
$ws.Range("I13").Value = 0.2478383111533466
$ws.Range("J13").Value = 0.5541963385427369
$ws.Range("K13").Value = -0.05728305666909728
$ws.Range("B14").Value = 0.9871738597754134
$ws.Range("C14").Value = 0.4122000866690486
$ws.Range("D14").Value = 0.3051903340877904
$ws.Range("E14").Value = 0.6592080140502106
$ws.Range("F14").Value = 0.1487291986403278
$ws.Range("G14").Value = 0.3684124016882561
$ws.Range("H14").Value = 0.2289442202646768
$ws.Range("I14").Value = 0.5393323377276911
$ws.Range("J14").Value = -0.085994767298468
$ws.Range("K14").Value = 0.6071338394308724
$ws.Range("B15").Value = 0.8628949586592991
$ws.Range("C15").Value = 0.3523010363001488
$ws.Range("D15").Value = 0.4189247832594023
$ws.Range("E15").Value = 0.1846742797061906
$ws.Range("F15").Value = 0.3663616442486632
$ws.Range("G15").Value = 0.1523759808286466
$ws.Range("H15").Value = 0.5119329433524077
$ws.Range("I15").Value = -0.1122838174961934
$ws.Range("J15").Value = 0.5637365487175399
$ws.Range("B16").Value = 0.6646266232236873
$ws.Range("C16").Value = 0.5522135229949265
$ws.Range("D16").Value = 0.005595857889999001
$ws.Range("E16").Value = 0.3945558830192304
$ws.Range("F16").Value = 0.1884012920210401
$ws.Range("G16").Value = 0.4743913731481941
$ws.Range("H16").Value = -0.1178261565835546
$ws.Range("I16").Value = 0.5679118776562884
$ws.Range("B17").Value = 0.7878040141027678
$ws.Range("C17").Value = 0.09027459876430857
$ws.Range("D17").Value = 0.2489555163855133
$ws.Range("E17").Value = 0.1983550348802827
$ws.Range("F17").Value = 0.489756542847739
$ws.Range("G17").Value = -0.1566744915409777
$ws.Range("H17").Value = 0.5520660686291026
$ws.Range("B18").Value = 0.4012987852456914
$ws.Range("C18").Value = 0.3660442907967085
$ws.Range("D18").Value = 0.03589089343071589
$ws.Range("E18").Value = 0.5187154933129405
$ws.Range("F18").Value = -0.1209316190860741
$ws.Range("G18").Value = 0.5207384222454754
$ws.Range("B19").Value = 0.6128694682008229
$ws.Range("C19").Value = 0.05354676096860539
$ws.Range("D19").Value = 0.4240932542019461
$ws.Range("E19").Value = -0.088714230837594
$ws.Range("F19").Value = 0.5331265480731927
$ws.Range("B20").Value = 0.2932231203848173
$ws.Range("C20").Value = 0.5091174976711597
$ws.Range("D20").Value = -0.2051988091956081
$ws.Range("E20").Value = 0.5482318197250452
$ws.Range("B21").Value = 0.6739775747052469
$ws.Range("C21").Value = -0.1916615369476919
$ws.Range("D21").Value = 0.4852589007350822
$ws.Range("B22").Value = 0.06218750491771613
$ws.Range("C22").Value = 0.5845769509171186
$ws.Range("B23").Value = 0.6286366421565677

$ws.Range("K15").ClearContents()
$ws.Range("J16").ClearContents()
$ws.Range("I17").ClearContents()
$ws.Range("H18").ClearContents()
$ws.Range("G19").ClearContents()
$ws.Range("F20").ClearContents()
$ws.Range("E21").ClearContents()
$ws.Range("D22").ClearContents()
$ws.Range("C23").ClearContents()
$ws.Range("B24").ClearContents()
